# Rename header columns to reflect clearer naming
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "season_ending_year_x"
$ws.Range("O1").Value = "season_ending_year_y"

# Fill in previously-missing birth_year values for the two rows
$ws.Range("Q2").Value = 1989
$ws.Range("Q3").Value = 1998

# Add new calendar_year column at the end of the table, matching the
# existing header formatting (bold, bordered, centered)
$ws.Range("AY1").Value = "calendar_year"
$ws.Range("AX1").Copy()
$ws.Range("AY1").PasteSpecial(-4122)
$ws.Range("AY1").Value = "calendar_year"

$ws.Range("AY2").Value = 2024
$ws.Range("AY3").Value = 2023
